$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns (AC, AD, AE)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AB1) to the new headers
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill the team record values for each data row (rows 2 through 51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 29).Value = 68   # AC -> Wins
    $ws.Cells.Item($r, 30).Value = 94   # AD -> Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE -> Ties
}
